# Update "想去人数" (F column) values on the 展览 sheet and the 全部类型 sheet
# to reflect newly scraped counts (gh-pages output update).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 1439
$ws1.Range("F15").Value = 1356
$ws1.Range("F26").Value = 5872
$ws1.Range("F31").Value = 14518
$ws1.Range("F36").Value = 9061
$ws1.Range("F37").Value = 624
$ws1.Range("F38").Value = 4211
$ws1.Range("F39").Value = 143
$ws1.Range("F40").Value = 361

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 1439
$ws4.Range("F15").Value = 1356
$ws4.Range("F29").Value = 5872
$ws4.Range("F34").Value = 14518
$ws4.Range("F39").Value = 9061
$ws4.Range("F40").Value = 624
$ws4.Range("F41").Value = 4211
$ws4.Range("F42").Value = 143
$ws4.Range("F43").Value = 361
